$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-row values for columns C..O (rows 2..25), per case "380 kV done".
$data = @(
    @(0.1891520658211192,0.03931962701445002,0.1208682022721383,0.9328533913881643,0.002453768097947433,0,0.9373883885711152,0,0.9945855932279528,0.1487661459673681,0.3441175607624842,1.344281385537869,3.299978809398993),
    @(0.1870229643982242,0.03785575005071706,0.1206388702552061,0.9312951537723677,0.002456560600396716,0,0.9386655667246231,0,0.8910802721925393,0.149282967225389,0.3219362935698129,1.352335726754973,3.309875331244683),
    @(0.1857942192809929,0.03694852640041546,0.1205543960241826,0.9309331599552593,0.00245836768979845,0,0.9399985151810881,0,0.8275414444452167,0.1496738887258466,0.3084097340383778,1.357769693075021,3.317962324230137),
    @(0.185313306601202,0.03657673512703496,0.1205341549926278,0.930935135522077,0.002459127417911109,0,0.940679658659576,0,0.8016539523775918,0.1498517011964893,0.3029212375434724,1.360107106320804,3.321763091946536),
    @(0.1852346498189803,0.03651487385069174,0.12053165108785,0.9309444904491784,0.002459254980947765,0,0.9408010932795676,0,0.7973557093542638,0.1498823448184048,0.3020113176204049,1.360502667783869,3.322424713759574),
    @(0.1857876532261571,0.03694352072856333,0.120554065598256,0.9309325814076459,0.002458377841268147,0,0.9400071427858165,0,0.827192293902101,0.1496762118206547,0.3083356178937038,1.357800717880274,3.318011537369784),
    @(0.1884016893764482,0.03881663744353148,0.1207774455372039,0.9321926217295982,0.00245471179993548,0,0.9377148430629418,0,0.958895016010274,0.1489290740637657,0.3364503921832807,1.346957222208943,3.302973665008153),
    @(0.1941485760460893,0.04242240850384604,0.1216619319857024,0.939387484449739,0.002448253342972201,0,0.9375764390816315,0,1.217215515216424,0.1480478085402908,0.392307752522548,1.329562553180118,3.289453847146433),
    @(0.1987467561009595,0.04502966710373357,0.1225834068925593,0.947562164113549,0.002443949302730236,0,0.9401361360117306,0,1.406975846295836,0.1477564030671132,0.4337757134402622,1.319132572433233,3.289284894229581),
    @(0.2009197328051044,0.04620652090012101,0.1230615011737051,0.9519103036338095,0.002442086097789071,0,0.9418797720017693,0,1.493285298549154,0.1477011895168374,0.4527315547516082,1.314896194797669,3.291334726476464),
    @(0.2017542105537444,0.04665082310965829,0.1232510019047481,0.9536474629267815,0.002441394100595556,0,0.9426234052761657,0,1.525965125499226,0.1476914047239575,0.4599225591839797,1.313364941904808,3.29241720406705),
    @(0.2015739751143713,0.04655519495311466,0.1232098136438253,0.9532693029452446,0.002441542532569271,0,0.9424595422183302,0,1.518927131986288,0.1476930173258211,0.4583732825998297,1.313691481006003,3.292170445653113),
    @(0.2009881533026601,0.04624310103697127,0.1230769221089503,0.9520514041752079,0.00244202889535087,0,0.9419392802033428,0,1.495973971542412,0.1477001615823355,0.4533229078258998,1.314768755900218,3.29141764235257),
    @(0.2006308317236858,0.04605175840948306,0.1229966231020505,0.9513172095104352,0.00244232857039799,0,0.9416314620672708,0,1.481913965614012,0.1477059862623875,0.4502310684511741,1.315438117507341,3.290996424359093),
    @(0.1986063762954302,0.04495256976475304,0.1225533465027659,0.9472906774999501,0.002444072968185951,0,0.9400338462303495,0,1.401334899619656,0.147761567914408,0.4325387213742857,1.319419647242356,3.289193764535355),
    @(0.1973852009508903,0.04427587962678814,0.1222964888066862,0.9449818165253276,0.002445167315991986,0,0.9392021587253581,0,1.351897553661274,0.1478154767609325,0.421708302803431,1.32199228219968,3.288632897420172),
    @(0.196690463839559,0.04388580070584425,0.1221542968759763,0.9437130603404853,0.002445805675460732,0,0.9387783088202042,0,1.323461336029879,0.1478537640399651,0.4154876006909447,1.323519841510091,3.288510449629513),
    @(0.1964565537432037,0.04375357890315712,0.1221071060732868,0.9432936526620495,0.002446023347053572,0,0.9386441614879999,0,1.313833174909917,0.147867977835773,0.4133828779081909,1.324045268416832,3.288503353278202),
    @(0.1975144057993816,0.04434800410694351,0.1223232578801721,0.9452214669076398,0.002445049897984495,0,0.9392850507313355,0,1.357160376687887,0.1478089845793846,0.4228603249083704,1.321713469624427,3.28867188493831),
    @(0.2011599085507072,0.04633480732501027,0.1231157261742837,0.9524066705490668,0.002441885671352171,0,0.9420898309860632,0,1.502715985228861,0.1476977612494039,0.4548059798287767,1.314450354638069,3.29163044391376),
    @(0.203610148042074,0.04762543707266786,0.1236829360141201,0.9576308259235304,0.002439896664733233,0,0.9444088106918613,0,1.597822674750375,0.1476899063265975,0.4757589745434032,1.31012877143425,3.295349297295871),
    @(0.2022962362895697,0.04693733077021278,0.1233757005609171,0.9547942315540752,0.002440951027388922,0,0.9431266467365731,0,1.547065064196147,0.1476881658846239,0.4645692574961728,1.312396405817907,3.293200979983339),
    @(0.1974559694415774,0.04431539986991595,0.1223111385146751,0.9451129383205057,0.002445102953909524,0,0.9392474061233216,0,1.354781097963723,0.1478118969738631,0.4223394774655276,1.321839369649787,3.288653635713018),
    @(0.192527691555938,0.04145425200770347,0.121374914274238,0.9369344444068659,0.002449922766619953,0,0.9371469015441463,0,1.14733357697861,0.1482236949568403,0.392307752522548,1.333854979548228,3.291398715648853)
)

$startRow = 2
$startCol = 3  # column C

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $startCol + $j).Value = $rowValues[$j]
    }
}
